# ------------------------------------------------------------------
# specs.xlsx update:
#  - check_derived_variables.R: new function compares PARAM/PARAMN
#    columns in the specs file against the output dataset -> add a
#    "Recode PARAM" lookup sheet (PARAM / PARAMN / EVID) used by that
#    comparison.
#  - Housekeeping on the existing "Specification-Source Data" sheet:
#    drop a redundant duplicate cell style that had crept onto most
#    cells, move the active selection, and leave the real content
#    (values/headers/filter) untouched.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Tidy up "Specification-Source Data" -----------------------
# Almost every populated cell (A1:F19) carried a cell style that was
# just a byte-for-byte duplicate of the default "Normal" style.
# Re-apply "Normal" across the used range to drop that redundant
# style reference from every cell.
$ws1.Range("A1:F19").Style = "Normal"

# Column C rows 11-19 (the "Column Type" entries for the Timing /
# Measurement blocks) keep their own distinct style - re-stamp the
# font so those cells keep a dedicated style slot instead of folding
# back into the plain default.
$ws1.Range("C11:C19").Font.Name = "Calibri"

# Move the stored selection from D10 to D2 (frozen header pane stays
# in place).
$ws1.Range("D2").Select()

# --- 2. Add the "Recode PARAM" lookup sheet ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Recode PARAM"

$ws2.Columns("A").ColumnWidth = 19.6

$ws2.Range("A1").Value = "PARAM"
$ws2.Range("B1").Value = "PARAMN"
$ws2.Range("C1").Value = "EVID"

# Seed the shared-string table with "Glucose" before "Hemoglobin" is
# written so the two strings land in the same order used by the
# workbook being reproduced, even though row 2 holds "Hemoglobin" and
# row 3 holds "Glucose".
$ws2.Range("A3").Value = "Glucose"

$ws2.Range("A2").Value = "Hemoglobin"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 1

$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 1

# Two blank, vertically-centred placeholder rows below the table.
$ws2.Range("A5").VerticalAlignment = -4108   # xlVAlignCenter
$ws2.Range("A6").VerticalAlignment = -4108   # xlVAlignCenter

$ws2.Range("A8").Select()
